$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2").Value = -38915.87
$ws.Range("D2").Value = -42885.77
$ws.Range("E2").Value = -38293.06
$ws.Range("F2").Value = -31954.07
$ws.Range("G2").Value = -28973.95
$ws.Range("H2").Value = -181022.72
$ws.Range("C3").Value = 57327.05
$ws.Range("D3").Value = 69089.91
$ws.Range("E3").Value = 46505.97
$ws.Range("F3").Value = 68042.99
$ws.Range("G3").Value = 97647.47
$ws.Range("H3").Value = 338613.39
$ws.Range("C4").Value = -13666.72
$ws.Range("D4").Value = -18150.04
$ws.Range("E4").Value = -17989.96
$ws.Range("F4").Value = -23666.17
$ws.Range("G4").Value = -17427.45
$ws.Range("H4").Value = -90900.34
$ws.Range("C5").Value = -22686.33
$ws.Range("D5").Value = -12307.34
$ws.Range("E5").Value = -20829.45
$ws.Range("F5").Value = -15279.92
$ws.Range("G5").Value = -17002.62
$ws.Range("H5").Value = -88105.66
$ws.Range("C6").Value = 452802.19
$ws.Range("D6").Value = 439948.04
$ws.Range("E6").Value = 521104.54
$ws.Range("F6").Value = 400785.61
$ws.Range("G6").Value = 439761.29
$ws.Range("H6").Value = 2254401.67
$ws.Range("C7").Value = -11866.61
$ws.Range("D7").Value = -20344.57
$ws.Range("E7").Value = -15438.74
$ws.Range("F7").Value = -15008.15
$ws.Range("G7").Value = -18061.25
$ws.Range("H7").Value = -80719.32
$ws.Range("C8").Value = -17838.35
$ws.Range("D8").Value = -22945.35
$ws.Range("E8").Value = -16119.28
$ws.Range("F8").Value = -19638.93
$ws.Range("G8").Value = -20240.81
$ws.Range("H8").Value = -96782.72
$ws.Range("C9").Value = -24322.07
$ws.Range("D9").Value = -17363.31
$ws.Range("E9").Value = -15047.77
$ws.Range("F9").Value = -18153.58
$ws.Range("G9").Value = -18839.59
$ws.Range("H9").Value = -93726.32
$ws.Range("C10").Value = -11891.18
$ws.Range("D10").Value = -20277.56
$ws.Range("E10").Value = -16132.08
$ws.Range("F10").Value = -18620.19
$ws.Range("G10").Value = -13311.51
$ws.Range("H10").Value = -80232.52
$ws.Range("C11").Value = -34389.91
$ws.Range("D11").Value = -32041.54
$ws.Range("E11").Value = -37326.52
$ws.Range("F11").Value = -56355.78
$ws.Range("G11").Value = -30370.45
$ws.Range("H11").Value = -190484.2
$ws.Range("C12").Value = -23546.99
$ws.Range("D12").Value = -12433.3
$ws.Range("E12").Value = -16701.58
$ws.Range("F12").Value = -17088.3
$ws.Range("G12").Value = -19492.48
$ws.Range("H12").Value = -89262.65
$ws.Range("C13").Value = -15592
$ws.Range("D13").Value = -18102.8
$ws.Range("E13").Value = -12464.41
$ws.Range("F13").Value = -15161.04
$ws.Range("G13").Value = -12295.02
$ws.Range("H13").Value = -73615.27
$ws.Range("C14").Value = -12972.63
$ws.Range("D14").Value = -11891.43
$ws.Range("E14").Value = -15180.02
$ws.Range("F14").Value = -14769.04
$ws.Range("G14").Value = -12670.83
$ws.Range("H14").Value = -67483.95
$ws.Range("C15").Value = 196582.08
$ws.Range("D15").Value = 188143.44
$ws.Range("E15").Value = 185434.74
$ws.Range("F15").Value = 205682.74
$ws.Range("G15").Value = 189920.85
$ws.Range("H15").Value = 965763.85
$ws.Range("C16").Value = -196582.08
$ws.Range("D16").Value = -188143.44
$ws.Range("E16").Value = -185434.74
$ws.Range("F16").Value = -205682.74
$ws.Range("G16").Value = -189920.85
$ws.Range("H16").Value = -965763.85
$ws.Range("C17").Value = -18135.67
$ws.Range("D17").Value = -16482.58
$ws.Range("E17").Value = -15629.56
$ws.Range("F17").Value = -19907.67
$ws.Range("G17").Value = -12846.96
$ws.Range("H17").Value = -83002.44
$ws.Range("C18").Value = 39584.34
$ws.Range("D18").Value = 40244.27
$ws.Range("E18").Value = 35541.82
$ws.Range("F18").Value = 37747.69
$ws.Range("G18").Value = 56504.28
$ws.Range("H18").Value = 209622.4
$ws.Range("C19").Value = 42131.76
$ws.Range("D19").Value = 36010.94
$ws.Range("E19").Value = 40646.08
$ws.Range("F19").Value = 30833.3
$ws.Range("G19").Value = 54080.06
$ws.Range("H19").Value = 203702.14
$ws.Range("C20").Value = 346021.01
$ws.Range("D20").Value = 340067.57
$ws.Range("E20").Value = 406645.98
$ws.Range("F20").Value = 271806.75
$ws.Range("G20").Value = 426460.18
$ws.Range("H20").Value = 1791001.49
